$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: column G "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-09-02 04:16:43"
$wsOverview.Range("G5").Value = "2016-09-02 04:16:43"

# zh-cn sheet: column E "Priority" ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: column H "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-02 04:16:39"
$wsZhCn.Range("H5").Value = "2016-09-02 04:16:39"

# zh-cn sheet: column K "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-09-02 04:16:56"
$wsZhCn.Range("K5").Value = "2016-09-02 04:16:56"

# de-de sheet: column H "Correspond Handoff Datetime"
$wsDeDe.Range("H2").Value = "2016-09-02 04:16:43"
$wsDeDe.Range("H5").Value = "2016-09-02 04:16:43"

# de-de sheet: column K "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-09-02 04:17:08"
$wsDeDe.Range("K5").Value = "2016-09-02 04:17:08"
